$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cryptocurrency price/volume data to reflect the latest GitHub Actions scrape.
# Each cell holds its value as literal text (matching the original inlineStr cells),
# so NumberFormat is forced to Text ("@") before the write to stop Excel from
# auto-coercing the numeric-looking / percentage-looking strings into numbers.

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "326.70"
$ws.Range("E2").Value = "-2.59%"

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "44.59"
$ws.Range("E3").Value = "1.85%"

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "5.609"
$ws.Range("E4").Value = "-2.81%"

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08056"
$ws.Range("E5").Value = "-3.24%"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.38%"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-4.85%"

$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9454"
$ws.Range("E9").Value = "0.23%"

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1166"
$ws.Range("E10").Value = "-6.28%"

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1860"
$ws.Range("E11").Value = "-5.08%"

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09878"
$ws.Range("E12").Value = "-0.81%"

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04245"
$ws.Range("E13").Value = "-7.04%"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.14%"

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001288"
$ws.Range("E15").Value = "-0.43%"

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04213"
$ws.Range("E16").Value = "-4.94%"

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005863"
$ws.Range("E17").Value = "-1.02%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.600"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.15%"

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "8.428"
$ws.Range("E20").Value = "-4.02%"

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1371"
$ws.Range("E21").Value = "0.59%"

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2612"
$ws.Range("E22").Value = "-3.07%"

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001245"
$ws.Range("E23").Value = "-1.47%"

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004465"
$ws.Range("E24").Value = "2.67%"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.42%"

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003994"
$ws.Range("E26").Value = "-0.06%"

$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02640"
$ws.Range("E38").Value = "-5.86%"

$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05473"
$ws.Range("E39").Value = "-6.15%"

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007697"
$ws.Range("E40").Value = "-3.12%"

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1397"
$ws.Range("E41").Value = "-2.27%"

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007014"
$ws.Range("E42").Value = "-21.89%"

$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002011"
$ws.Range("E43").Value = "-7.49%"

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008857"
$ws.Range("E44").Value = "-15.19%"

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00007136"
$ws.Range("E45").Value = "1.74%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.07%"

$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "0.003681"
$ws.Range("E47").Value = "15.34%"

$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002273"
$ws.Range("E48").Value = "-0.06%"

$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002103"
$ws.Range("E49").Value = "-0.07%"

$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002003"
$ws.Range("E50").Value = "-0.07%"
